# Update the "timestamp" column (O) for all data rows (2..64) from the
# old scrape timestamp to the new one, matching the commit's re-run of
# the crawler on 2022-09-07 at 21:02:19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2022-09-07 07:14:51"
$newTimestamp = "2022-09-07 21:02:19"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 15).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 64 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 15) # column O
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
